$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2095588235294118
$ws.Range("C2").Value = 0.5588235294117647
$ws.Range("J2").Value = 0.02205882352941177
$ws.Range("P2").Value = 0.1360294117647059
$ws.Range("S2").Value = 0.07352941176470588
$ws.Range("B3").Value = 0.0124223602484472
$ws.Range("C3").Value = 0.03726708074534162
$ws.Range("J3").Value = 0.01863354037267081
$ws.Range("P3").Value = 0.7639751552795031
$ws.Range("S3").Value = 0.1677018633540373
$ws.Range("J4").Value = 0.06060606060606061
$ws.Range("P4").Value = 0.696969696969697
$ws.Range("S4").Value = 0.2424242424242424
$ws.Range("B6").Value = 0.05288461538461538
$ws.Range("D6").Value = 0.004807692307692308
$ws.Range("E6").Value = 0.004807692307692308
$ws.Range("F6").Value = 0.0625
$ws.Range("J6").Value = 0.25
$ws.Range("O6").Value = 0.04326923076923077
$ws.Range("Q6").Value = 0.1538461538461539
$ws.Range("R6").Value = 0.07211538461538461
$ws.Range("S6").Value = 0.3557692307692308
$ws.Range("B7").Value = 0.1069182389937107
$ws.Range("D7").Value = 0.01257861635220126
$ws.Range("F7").Value = 0.03773584905660377
$ws.Range("J7").Value = 0.169811320754717
$ws.Range("O7").Value = 0.01886792452830189
$ws.Range("Q7").Value = 0.1886792452830189
$ws.Range("R7").Value = 0.09433962264150944
$ws.Range("S7").Value = 0.3710691823899371
$ws.Range("B8").Value = 0.08398950131233596
$ws.Range("D8").Value = 0.01837270341207349
$ws.Range("F8").Value = 0.05774278215223097
$ws.Range("J8").Value = 0.1286089238845144
$ws.Range("O8").Value = 0.02362204724409449
$ws.Range("Q8").Value = 0.1758530183727034
$ws.Range("R8").Value = 0.1076115485564304
$ws.Range("S8").Value = 0.4041994750656168
$ws.Range("B9").Value = 0.0947867298578199
$ws.Range("D9").Value = 0.01895734597156398
$ws.Range("F9").Value = 0.07582938388625593
$ws.Range("J9").Value = 0.1469194312796208
$ws.Range("O9").Value = 0.04265402843601896
$ws.Range("Q9").Value = 0.1658767772511848
$ws.Range("R9").Value = 0.0947867298578199
$ws.Range("S9").Value = 0.3601895734597156
$ws.Range("B10").Value = 0.1117696867061812
$ws.Range("D10").Value = 0.01947502116850127
$ws.Range("F10").Value = 0.07197290431837426
$ws.Range("J10").Value = 0.1160033869602032
$ws.Range("O10").Value = 0.01524132091447926
$ws.Range("Q10").Value = 0.1947502116850127
$ws.Range("R10").Value = 0.0821337849280271
$ws.Range("S10").Value = 0.388653683319221
$ws.Range("G11").Value = 0.1376518218623482
$ws.Range("J11").Value = 0.08502024291497975
$ws.Range("K11").Value = 0.1983805668016194
$ws.Range("L11").Value = 0.5627530364372469
$ws.Range("S11").Value = 0.01619433198380567
$ws.Range("G12").Value = 0.7887323943661971
$ws.Range("J12").Value = 0.1690140845070423
$ws.Range("K12").Value = 0.01408450704225352
$ws.Range("S12").Value = 0.02816901408450704
$ws.Range("G13").Value = 0.5945945945945946
$ws.Range("J13").Value = 0.3243243243243243
$ws.Range("S13").Value = 0.08108108108108109
$ws.Range("S14").Value = 1
$ws.Range("F15").Value = 0.0186046511627907
$ws.Range("H15").Value = 0.1441860465116279
$ws.Range("I15").Value = 0.06511627906976744
$ws.Range("J15").Value = 0.3116279069767442
$ws.Range("K15").Value = 0.03255813953488372
$ws.Range("M15").Value = 0.0186046511627907
$ws.Range("O15").Value = 0.06976744186046512
$ws.Range("S15").Value = 0.3395348837209302
$ws.Range("F16").Value = 0.01685393258426966
$ws.Range("H16").Value = 0.1292134831460674
$ws.Range("I16").Value = 0.1348314606741573
$ws.Range("J16").Value = 0.4157303370786517
$ws.Range("K16").Value = 0.06179775280898876
$ws.Range("M16").Value = 0.02247191011235955
$ws.Range("N16").Value = 0.005617977528089887
$ws.Range("O16").Value = 0.05617977528089887
$ws.Range("S16").Value = 0.1573033707865168
$ws.Range("F17").Value = 0.0103359173126615
$ws.Range("H17").Value = 0.1679586563307494
$ws.Range("I17").Value = 0.1266149870801034
$ws.Range("J17").Value = 0.4315245478036176
$ws.Range("K17").Value = 0.07235142118863049
$ws.Range("M17").Value = 0.01808785529715762
$ws.Range("N17").Value = 0.002583979328165375
$ws.Range("O17").Value = 0.041343669250646
$ws.Range("S17").Value = 0.1291989664082687
$ws.Range("F18").Value = 0.02127659574468085
$ws.Range("H18").Value = 0.1595744680851064
$ws.Range("I18").Value = 0.101063829787234
$ws.Range("J18").Value = 0.4148936170212766
$ws.Range("K18").Value = 0.1117021276595745
$ws.Range("M18").Value = 0.01595744680851064
$ws.Range("O18").Value = 0.06914893617021277
$ws.Range("S18").Value = 0.1063829787234043
$ws.Range("F19").Value = 0.02358887952822241
$ws.Range("H19").Value = 0.1946082561078349
$ws.Range("I19").Value = 0.08845829823083404
$ws.Range("J19").Value = 0.3740522325189554
$ws.Range("K19").Value = 0.1078348778433024
$ws.Range("M19").Value = 0.0160067396798652
$ws.Range("O19").Value = 0.07497893850042123
$ws.Range("S19").Value = 0.1204717775905644
